$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text (numeric-looking) cell value while preserving default (unstyled) formatting
function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "238.36"

Set-TextCell "D3" "21.62"

Set-TextCell "D4" "5.361"

Set-TextCell "D5" "0.05553"

Set-TextCell "D6" "3.364"

Set-TextCell "D7" "6.454"

Set-TextCell "D8" "0.8044"

Set-TextCell "D9" "1.047"

Set-TextCell "D11" "0.07289"

Set-TextCell "D12" "0.03259"

Set-TextCell "D13" "0.02870"

Set-TextCell "D14" "0.09235"

Set-TextCell "D15" "0.001659"

Set-TextCell "D16" "3.248"

Set-TextCell "D17" "0.04752"

Set-TextCell "D18" "0.0005701"
$ws.Range("E18").Value = "17OneONE"

Set-TextCell "D19" "0.006251"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell "D20" "0.001051"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell "D21" "0.003774"
$ws.Range("E21").Value = "20HotbitTokenHTB"

Set-TextCell "D22" "0.0001498"

Set-TextCell "D23" "0.0004175"

Set-TextCell "D24" "3.944"

Set-TextCell "D25" "2.201"

Set-TextCell "D40" "0.04152"

Set-TextCell "D41" "0.006997"

Set-TextCell "D42" "0.003495"

Set-TextCell "D43" "0.1038"

Set-TextCell "D44" "0.008799"

Set-TextCell "D45" "0.00005433"

Set-TextCell "D46" "0.00000000749"

Set-TextCell "D47" "0.6790"

Set-TextCell "D48" "0.03147"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

Set-TextCell "D49" "0.00002097"

Set-TextCell "D50" "0.01009"

Write-Output "applied symbol list update"